# "Generate Report for Archive"
#
# 1. Status text "Ready for handoff" -> "In Translation" on every sheet
#    that carries it (Overview columns E/F, and the "Status" column (C)
#    of each per-locale sheet).
# 2. Narrow the "zh-cn"/"de-de" status columns (Overview!E:F and the
#    "Status" column on each locale sheet) from ~17.22 chars to ~13.41
#    chars, matching the narrower report column width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($oldStatus, $newStatus)
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = 12.5
